$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column at H (shifts the old "2013-Oct-xx" H:I block to J:K,
#        and moves the H2:I2 merged header cell along with it). ---
$ws.Range("H1:I1").EntireColumn.Insert()

# --- 2. New "2013-Dec-xx" column header + labels (values intentionally blank,
#        matching the source data which only lists table names here). New
#        vocabulary is introduced in alphabetical order (degrees, milestones,
#        objects, offices, relationships), matching how the source list was
#        authored. ---
$ws.Range("H2").Value = "2013-Dec-xx"
$ws.Range("H26").Value = "degrees"
$ws.Range("H27").Value = "milestones"
$ws.Range("H23").Value = "objects"
$ws.Range("H24").Value = "offices"
$ws.Range("H25").Value = "relationships"
$ws.Range("H3").Value = "funding_rounds"
$ws.Range("H4").Value = "acquisitions"
$ws.Range("H6").Value = "investments"
$ws.Range("H7").Value = "ipos"
$ws.Range("H9").Value = "people"
$ws.Range("H10").Value = "funds"

# --- 3. Bring "people" and "funds" to the top of the 2016-Sep-09 (D:E) list,
#        keeping the remaining rows in their original relative order. ---
$ws.Range("D9").Value = "people"
$ws.Range("E9").Value = 545451
$ws.Range("D10").Value = "funds"
$ws.Range("E10").Value = 4954

$ws.Range("D16").Value = "people_descriptions"
$ws.Range("E16").Value = 290227
$ws.Range("D17").Value = "competitors"
$ws.Range("E17").Value = 502358
$ws.Range("D18").Value = "event_relationships"
$ws.Range("E18").Value = 109569
$ws.Range("D19").Value = "events"
$ws.Range("E19").Value = 7075
$ws.Range("D20").Value = "customers"
$ws.Range("E20").Value = 251337
$ws.Range("D21").Value = "investment_partners"
$ws.Range("E21").Value = 40994

$ws.Range("D22:E23").Clear()

# --- 4. "additions" row moves up from 24 to 22. ---
$ws.Range("F22").Value = "additions"
$ws.Range("G22").Value = 2213
$ws.Range("F24:G24").Clear()

# --- 5. New merged header cell for the inserted "2013-Dec-xx" column. ---
$ws.Range("H2:I2").Merge()

# --- 6. Column width for the new H column (bestFit-style autosize to its
#        longest label, "funding_rounds"). ---
$ws.Columns.Item(8).ColumnWidth = 13.6666666

# --- 7. Selection moves to where the user was last working. ---
$ws.Range("H17").Select()
